$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Price column (D) to text format so numeric-looking
# strings like "68.646.31" are not coerced into numbers/dates by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '68.646.31'
$ws.Range("D3").Value = '3.718.52'
$ws.Range("D5").Value = '600.88'
$ws.Range("D6").Value = '166.88'
$ws.Range("D7").Value = '3.715.41'
$ws.Range("D9").Value = '0.536'
$ws.Range("D10").Value = '0.163'
$ws.Range("D11").Value = '6.33'
$ws.Range("D13").Value = '37.96'
$ws.Range("D14").Value = '0.0000244'
$ws.Range("D15").Value = '4.339.79'
$ws.Range("D16").Value = '3.717.06'
$ws.Range("D17").Value = '68.603.21'
$ws.Range("D18").Value = '7.27'
$ws.Range("D20").Value = '17.05'
$ws.Range("D21").Value = '496.42'
$ws.Range("D22").Value = '10.39'
$ws.Range("D23").Value = '0.724'
$ws.Range("D24").Value = '85.27'
$ws.Range("D25").Value = '0.0000143'
$ws.Range("D26").Value = '2.30'
$ws.Range("D27").Value = '12.38'
$ws.Range("D28").Value = '10.14'
$ws.Range("D30").Value = '2.59'
$ws.Range("D31").Value = '2.96'
$ws.Range("D32").Value = '7.94'
$ws.Range("D33").Value = '31.39'
$ws.Range("D34").Value = '3.863.80'
$ws.Range("D36").Value = '3.650.73'
$ws.Range("D40").Value = '0.131'
$ws.Range("D42").Value = '435.23'
$ws.Range("D43").Value = '48.73'
$ws.Range("D44").Value = '1.97'
$ws.Range("D45").Value = '2.85'
$ws.Range("D46").Value = '8.50'
$ws.Range("D48").Value = '40.41'
$ws.Range("D49").Value = '141.42'
$ws.Range("D50").Value = '0.0352'
$ws.Range("D51").Value = '2.758.84'

# Restore the default (no explicit) style on the Price column so the
# workbook formatting matches the original.
$ws.Range("D2:D51").Style = "Normal"

# Volume(1h) column (E) values are plain text already (leading/trailing
# spaces + % sign), so they round-trip as text without extra handling.
$ws.Range("E2").Value = '  +1.31%  '
$ws.Range("E3").Value = '  -2.57%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("E5").Value = '  +0.07%  '
$ws.Range("E6").Value = '  -3.75%  '
$ws.Range("E7").Value = '  -2.62%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  +0.82%  '
$ws.Range("E10").Value = '  +2.67%  '
$ws.Range("E11").Value = '  +0.97%  '
$ws.Range("E12").Value = '  -1.30%  '
$ws.Range("E13").Value = '  -1.55%  '
$ws.Range("E14").Value = '  -0.30%  '
$ws.Range("E15").Value = '  -2.56%  '
$ws.Range("E16").Value = '  -2.69%  '
$ws.Range("E17").Value = '  +1.14%  '
$ws.Range("E18").Value = '  +0.25%  '
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("E20").Value = '  +0.51%  '
$ws.Range("E21").Value = '  +0.81%  '
$ws.Range("E22").Value = '  +11.30%  '
$ws.Range("E23").Value = '  -2.96%  '
$ws.Range("E24").Value = '  -0.79%  '
$ws.Range("E25").Value = '  -3.26%  '
$ws.Range("E26").Value = '  -2.65%  '
$ws.Range("E27").Value = '  +0.70%  '
$ws.Range("E28").Value = '  -1.36%  '
$ws.Range("E29").Value = '  -0.18%  '
$ws.Range("E30").Value = '  +5.71%  '
$ws.Range("E31").Value = '  -0.73%  '
$ws.Range("E32").Value = '  +1.23%  '
$ws.Range("E33").Value = '  -6.21%  '
$ws.Range("E34").Value = '  -2.38%  '
$ws.Range("E35").Value = '  -1.36%  '
$ws.Range("E36").Value = '  -2.74%  '
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("E38").Value = '  -0.53%  '
$ws.Range("E39").Value = '  -0.08%  '
$ws.Range("E40").Value = '  -1.17%  '
$ws.Range("E41").Value = '  -1.77%  '
$ws.Range("E42").Value = '  -5.02%  '
$ws.Range("E43").Value = '  -0.78%  '
$ws.Range("E44").Value = '  -1.90%  '
$ws.Range("E45").Value = '  -1.09%  '
$ws.Range("E46").Value = '  +0.74%  '
$ws.Range("E48").Value = '  -1.78%  '
$ws.Range("E49").Value = '  +1.38%  '
$ws.Range("E50").Value = '  +0.18%  '
$ws.Range("E51").Value = '  -3.45%  '
